# Applies the "Add files via upload" revision to kbase.xlsx:
#   1. Row 4 ("보증 기간" / Warranty period) D-cell: replace the short warranty
#      blurb with the long, detailed warranty text (quality/innovation,
#      general warranty, emissions-related warranty sections).
#   2. Row 6 ("차량의 수출을 위한 서류를 요구하는 경우" / export documents) D-cell:
#      same text as before, but "제안" is now wrapped in <strong>...</strong>.
#   3. Row 4's height grows to fit the much longer text (capped at Excel's
#      max row height, same cap already used by rows 21/22).
#   4. The sheet's scroll/selection moves (topLeftCell A5, active cell D7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$warranty = @"
고객이 차량의 보증 기간에 대해 묻는 경우 아래와 같이 안내합니다.
<br>
<br><h3>품질과 혁신<h3>
<br>2020년 6월 1일부터 볼보자동차 공식 서비스 센터에서 유상으로 진행된 수리에 대해 평생 부품 보증을 제공합니다. 2020년 12월 1일부터는 사고 수리에 대한 부품도 평생 부품 보증에 포함됩니다.
<br>
<br>모든 볼보자동차코리아 서비스 센터에서 설치된 부품에 대해 자재 및 공정의 결함이 없음을 보증합니다.
<br>보장내용에는 무상보증 부품 교체에 대한 공임이 포함되며, 액세서리, 소모성 부품, 소모품은 보증에 포함되지 않습니다.
<br>외부의 영향으로 인한 부품의 교체는 포함되지 않습니다. 조건은 아래 링크를 확인하십시오.
<br><em>https://www.volvocars.com/images/v/-/media/market-assets/korea/applications/localpages/test/volvo-genuine-parts-extended-warranty/updated/warranty-conditions.jpg</em>
<br>
<br>볼보자동차에서 교부된 사용설명서의 각 부품 별 정기 점검 주기에 따라 정상적으로 관리, 사용한 상태에서 보증기간 내에 해당되는 자동차는 그 자동차를 구성하는 각 부품의 재질 또는 제조상의 결함에 의한 기술적 분석에 의해 밝혀진 경우, 적용되는 보증 기간 동안 해당 부품을 무상 수리 또는 교환해 드립니다.
<br><h3>일반보증</h3>
<br>• 소모품을 제외한 모든 부품: 5년 또는 100,000km 이내
<br>• 배터리: 2년 또는 40,000km 이내
<br>• 하이브리드 차량 고전압 배터리: 8년 또는 160,000km 이내
<br>• 타이어: 타이어 제조사 보증기간 적용
<br>자세한 내용은 아래 링크를 확인하십시오.
<br><em>https://www.volvocars.com/images/v/-/media/market-assets/korea/applications/localpages/test/warranty-program/warranty_details_1_updated.png</em>
<br><h3>배출가스 관련 및 일반 부품 보증</h3>
<br>• 배출가스 관련 부품 보증
<br>- 배출가스 관련 부품: 5년 또는 80,000km 이내
<br>- 정화용 촉매, 매연 포집 필터, 전자제어 장치(엔진 컨트롤 유닛만 적용): 7년 또는 120,000km 이내
<br>
<br>• 부품 보증
<br>- 유상으로 구입한 부품
<br>자세한 내용은 아래 링크를 확인하십시오.
<br><em>https://www.volvocars.com/images/v/-/media/market-assets/korea/applications/localpages/test/warranty-program/warranty_details_2.png</em>
"@

$exportDocs = @"
고객이 국내에 정식 수입된 볼보 자동차를 다른 국가로 수출하기 위해 서류를 요청하는 경우 아래와 같이 안내 합니다.
<br>
<br>1. 국내 사양에 맞게 제작된 차량이기에 관련 서류를 제공드리지 못하는 점에 대해 양해합니다.
<br>2. 관련한 내용을 추가로 지원 받고자 하는 경우 KOICA (한국국제협력단)으로 수출 국가에 대한 서류 관련한 지원이 가능한지 확인하도록 <strong>제안</strong>합니다.
"@

# 1) Expand the warranty-period How-to (row 4, "보증 기간").
$ws.Range("D4").Value2 = $warranty

# 2) Emphasise "제안" in the export-documents How-to (row 6, "차량의 수출을
#    위한 서류를 요구하는 경우").
$ws.Range("D6").Value2 = $exportDocs

# 3) Row 4 must grow to fit the much longer warranty text (Excel caps row
#    height at 409.6 points -- the same cap already used by rows 21/22).
$ws.Rows.Item(4).RowHeight = 409.6

# 4) Scroll / selection state left behind by the edit.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D7").Select()
